# Upload new version with timestamp
# - Bumps the "م" (serial number) column A values in rows 7-105 by +2
#   (they now continue from a prior batch, so 1..99 becomes 3..101).
# - Updates the generated/export timestamp text cell from 11:58 PM to 11:59 PM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 7; $row -le 105; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value()
    $cell.Value = $current + 2
}

$ws.Range("A106").Value = "Friday, 15 August, 2025 11:59 PM"
